# Swap the order of recorder names in column G ("Recorded By") for the
# specified rows: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# and "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowsSystem = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,30,33,37,38,39,40,41,42,44,45,46,47,48,49,51,57,60,64,65,66,67,68,69,71,72,73,74,75,76,78,86,87,88,89,93,95,96,97,99,102,112,113,114,115,119,121,122,123,125,128,138,139,140,141,145,147,148,149,151,154)
$rowsAdmin = @(90,116,142)

foreach ($r in $rowsSystem) {
    $cell = $ws.Range("G$r")
    if ($cell.Text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

foreach ($r in $rowsAdmin) {
    $cell = $ws.Range("G$r")
    if ($cell.Text -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value = "admin@admin.com, dnasr281@gmail.com"
    }
}

$wb.Save()
